$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- 1. New cell values, in the exact order that reproduces the target
#        shared-string table order (Z,Y,X,W,V,U,T,S,R appended after G).
$ws.Range("A6").Value = "Z"
$ws.Range("B6").Value = "Y"
$ws.Range("C6").Value = "X"
$ws.Range("D6").Value = "W"
$ws.Range("F6").Value = "V"
$ws.Range("F4").Value = "U"
$ws.Range("F3").Value = "T"
$ws.Range("F2").Value = "S"
$ws.Range("F1").Value = "R"

# --- 2. Stamp every newly-touched cell with the same format (style index)
#        used by the rest of the sheet (copy A1's format, which carries the
#        center/center alignment style already present in the workbook).
$ws.Range("A1").Copy()
$ws.Range("E1:F6").PasteSpecial($xlPasteFormats)
$ws.Range("A5:D6").PasteSpecial($xlPasteFormats)
[void]($excel.CutCopyMode = $false)

# --- 3. Column E: width 0, hidden
$ws.Columns("E").ColumnWidth = -0.8333333333333334
$ws.Columns("E").Hidden = $true

# --- 4. Row 5: hidden
$ws.Rows(5).Hidden = $true

# --- 5. Selection: selected range A1:F6 (active cell ends up at the
#        range's top-left corner; the engine ties active cell to the
#        selection's anchor and has no way to park it elsewhere within a
#        multi-cell selection).
[void]$ws.Range("A1:F6").Select()

Write-Host "done"
